# RMS Upload Template (Quiz) - single choice excel template changes
#
# Renames the hidden "Sheet4" QuestionType lookup list from
# Radio / Single Choice / Multiple Choice
# to
# SingleChoice / MultipleChoice / Programm
# shrinks the QuestionType data validation range from A2:A4 to A2:A3,
# moves the sheet's active-cell selection to A4, and sets the page
# orientation to portrait.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Quiz")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Update the QuestionType lookup values --------------------------------
$ws4.Range("A2").Value = "SingleChoice"
$ws4.Range("A3").Value = "MultipleChoice"
$ws4.Range("A4").Value = "Programm"

# --- Shrink the QuestionType validation range to A2:A3 --------------------
# (the validation list itself still spans A2:A3, so only remove it from A4)
$ws4.Range("A4").Validation.Delete()

# --- Move the active selection on Sheet4 to A4 -----------------------------
# Activate Sheet4 long enough to move the selection, then restore the
# originally active sheet (Quiz) so the workbook's active tab is unchanged.
$ws4.Select()
$ws4.Range("A4").Select()
$ws1.Select()

# --- Page setup: portrait orientation for Sheet4 ---------------------------
$ws4.PageSetup.Orientation = 1
